# key_buildingBlock_pairs.xlsx — insert three blocks of 10 blank rows into
# Sheet1, pushing the existing "key/buildingBlock/pair" rows further down
# the sheet (old A11:C111 data block -> new A21:C141).
#
# Because we work from the bottom of the sheet upward, each Rows(...).Insert()
# call below can be expressed using the worksheet's *original* row numbers:
# rows below a later (higher-numbered) insertion point haven't shifted yet
# when we operate on an earlier (lower-numbered) one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 10 rows above original row 91 (old rows 91/101/111 -> 121/131/141).
$ws.Rows("91:100").Insert()

# 2) Insert 10 rows above original row 51 (old rows 51-84 -> 71-104).
$ws.Rows("51:60").Insert()

# 3) Insert 10 rows above original row 11 (old rows 11-45 -> 21-55).
$ws.Rows("11:20").Insert()

# Match the saved view/selection state: active cell B112, scrolled so row 94
# is at the top of the window.
$ws.Range("A94").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 94
$ws.Range("B112").Select() | Out-Null
